# Update the "acquisition datetime" (取得日時) timestamps in column A
# of the "ランサーズ" sheet for the existing data rows (2-10) to reflect
# the new run time: 2025-12-03 12:52:47 (JST).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-03 12:52:47"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
